# Rename the "Region" model/sheet to "Zone" (per commit message: renaming
# model Region to Zone) and restore the Excel UI state that comes along
# with a user opening that sheet and clicking a cell on it before saving
# (tabSelected moves from "Line" to "Zone", the active workbook tab index
# updates, and the last selected cell on the renamed sheet becomes D3).

$wb = $excel.ActiveWorkbook

$zone = $wb.Worksheets.Item("Region")
$zone.Name = "Zone"

# Make "Zone" the active sheet -> updates workbookView/@activeTab and
# moves sheetView/@tabSelected from the previously active "Line" sheet.
$zone.Activate()

# Leave the last selection on the sheet at D3 (was D4).
$zone.Range("D3").Select() | Out-Null
